$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (headers) - extend header formatting (bold + border + center/top align)
# from C1 into the new D1/E1 header cells, then set their text.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

$ws.Range("C1").Value = "OCDE"
$ws.Range("D1").Value = "Commission Européenne"
$ws.Range("E1").Value = "FMI"

# Row 2 (separators)
$ws.Range("D2").Value = "---"
$ws.Range("E2").Value = "---"

# Row 3
$ws.Range("A3").Value = "PIB (2024)"
$ws.Range("B3").Value = "1,1 %"
$ws.Range("C3").Value = "1,1 %"
$ws.Range("D3").Value = "0,7 %"
$ws.Range("E3").Value = "0,9 %"

# Row 4
$ws.Range("A4").Value = "PIB (2025)"
$ws.Range("B4").Value = "1,1 %"
$ws.Range("C4").Value = "1,2 %"
$ws.Range("D4").Value = "1,3 %"
$ws.Range("E4").Value = "1,3 %"

# Row 5
$ws.Range("A5").Value = "IPC (2024)"
$ws.Range("B5").Value = "2,5 %"
$ws.Range("C5").Value = "2,4 %"
$ws.Range("D5").Value = "2,5 %"
$ws.Range("E5").Value = "n.d."

# Row 6
$ws.Range("A6").Value = "IPC (2025)"
$ws.Range("B6").Value = "1,9 %"
$ws.Range("C6").Value = "1,9 %"
$ws.Range("D6").Value = "2,0 %"
$ws.Range("E6").Value = "n.d."
